# Auto-generated Excel COM-interop script applying scheduled-runner sheet updates
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
# across all class sheets, as produced by the scheduled price-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M21").Value = -2782
$ws.Range("H21").Value = 3250
$ws.Range("I21").Value = 3250
$ws.Range("K21").Value = 3250
$ws.Range("K23").Value = 3250
$ws.Range("H23").Value = 3250
$ws.Range("I23").Value = 3250
$ws.Range("M23").Value = -3016
$ws.Range("K28").Value = 887.25
$ws.Range("M28").Value = -402.25
$ws.Range("I28").Value = 887.25
$ws.Range("H28").Value = 989.8
$ws.Range("H33").Value = 125.09524
$ws.Range("I33").Value = 100.63636
$ws.Range("N33").Value = -610
$ws.Range("M33").Value = 128.36364
$ws.Range("L33").Value = 152
$ws.Range("J33").Value = 152
$ws.Range("K33").Value = 100.63636
$ws.Range("I80").Value = 1590
$ws.Range("K80").Value = 4770
$ws.Range("H80").Value = 2642.5
$ws.Range("M80").Value = -3772
$ws.Range("K83").Value = 14310
$ws.Range("I83").Value = 1590
$ws.Range("H83").Value = 2642.5
$ws.Range("M83").Value = -9318
$ws.Range("I94").Value = 4315.8335
$ws.Range("K94").Value = 4315.8335
$ws.Range("M94").Value = -3864.8335
$ws.Range("H94").Value = 4556.4287
$ws.Range("I100").Value = 1866.4445
$ws.Range("K100").Value = 1866.4445
$ws.Range("H100").Value = 1866.4445
$ws.Range("M100").Value = -1325.4445
$ws.Range("I107").Value = 177.28572
$ws.Range("M107").Value = 1742.71428
$ws.Range("K107").Value = 177.28572
$ws.Range("H107").Value = 177.28572
$ws.Range("I138").Value = 20002750
$ws.Range("J138").Value = 2984.0952
$ws.Range("H138").Value = 3849092.8
$ws.Range("K138").Value = 60008250
$ws.Range("L138").Value = 8952.285600000001
$ws.Range("N138").Value = -19232.2856
$ws.Range("M138").Value = -60003110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("I132").Value = 1190.3
$ws.Range("K132").Value = 3570.9
$ws.Range("H132").Value = 1669.4615
$ws.Range("M132").Value = -1040.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M134").Value = -2325.857400000001
$ws.Range("I134").Value = 1620.2858
$ws.Range("H134").Value = 2042.75
$ws.Range("K134").Value = 4860.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N32").Value = -5464.3335
$ws.Range("K32").Value = 3000.5
$ws.Range("I32").Value = 3000.5
$ws.Range("L32").Value = 4832.3335
$ws.Range("H32").Value = 4099.6
$ws.Range("J32").Value = 4832.3335
$ws.Range("M32").Value = -2684.5
$ws.Range("H68").Value = 49400
$ws.Range("L68").Value = 49400
$ws.Range("J68").Value = 49400
$ws.Range("N68").Value = -50898
$ws.Range("H71").Value = 49400
$ws.Range("J71").Value = 49400
$ws.Range("N71").Value = -155688
$ws.Range("L71").Value = 148200
$ws.Range("L141").Value = 695177
$ws.Range("N141").Value = -705537
$ws.Range("H141").Value = 695177
$ws.Range("J141").Value = 695177

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 423160.2
$ws.Range("K4").Value = 1269480.6
$ws.Range("J4").Value = 1000000000
$ws.Range("L4").Value = 3000000000
$ws.Range("N4").Value = -3000000224
$ws.Range("H4").Value = 167019300
$ws.Range("M4").Value = -1269368.6
$ws.Range("M50").ClearContents()
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M134").Value = 2370
$ws.Range("I134").Value = 900
$ws.Range("H134").Value = 1561
$ws.Range("K134").Value = 2700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M19").Value = -3528.875
$ws.Range("I19").Value = 3816.875
$ws.Range("K19").Value = 3816.875
$ws.Range("H19").Value = 3503.0908
$ws.Range("J57").Value = 50061
$ws.Range("L57").Value = 50061
$ws.Range("H57").Value = 50061
$ws.Range("N57").Value = -51701
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("I107").Value = 481.25
$ws.Range("M107").Value = 1438.75
$ws.Range("K107").Value = 481.25
$ws.Range("L107").Value = 625
$ws.Range("H107").Value = 529.1667
$ws.Range("N107").Value = -4465
$ws.Range("J107").Value = 625
$ws.Range("I113").Value = 1776.3334
$ws.Range("M113").Value = 393.6666
$ws.Range("K113").Value = 1776.3334
$ws.Range("H113").Value = 2095.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N21").Value = -17607.334
$ws.Range("H21").Value = 15444.5
$ws.Range("J21").Value = 17259.334
$ws.Range("L21").Value = 17259.334
$ws.Range("I40").Value = 1499.5
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1363.5
$ws.Range("N40").Value = -3272
$ws.Range("K40").Value = 1499.5
$ws.Range("H40").Value = 1999.6666
$ws.Range("I46").Value = 1787.75
$ws.Range("J46").Value = 2250
$ws.Range("M46").Value = -1599.75
$ws.Range("N46").Value = -2626
$ws.Range("L46").Value = 2250
$ws.Range("K46").Value = 1787.75
$ws.Range("H46").Value = 1880.2
$ws.Range("K61").Value = 1026.6666
$ws.Range("M61").Value = -824.6666
$ws.Range("I61").Value = 1026.6666
$ws.Range("J61").Value = 2124.25
$ws.Range("L61").Value = 2124.25
$ws.Range("H61").Value = 1364.3846
$ws.Range("N61").Value = -2528.25
$ws.Range("L113").Value = 2124.25
$ws.Range("J113").Value = 2124.25
$ws.Range("I113").Value = 1026.6666
$ws.Range("M113").Value = 1143.3334
$ws.Range("N113").Value = -6464.25
$ws.Range("K113").Value = 1026.6666
$ws.Range("H113").Value = 1364.3846

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J4").Value = 19250
$ws.Range("L4").Value = 19250
$ws.Range("N4").Value = -19476
$ws.Range("H4").Value = 18700.2
$ws.Range("K75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("H75").Value = 30000
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()
$ws.Range("K78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("H78").Value = 30000
$ws.Range("I122").Value = 1396
$ws.Range("M122").Value = -1738
$ws.Range("H122").Value = 1557.8182
$ws.Range("K122").Value = 4188
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
